# FedExShipments.xlsx - "Changes of File path in Address Book"
# The Address Book / tracking-number column (P) on Sheet1 gets refreshed
# with a new batch of FedEx tracking numbers for rows 2-26. The cell's
# underlying data type must stay text (not get auto-converted to a
# number) and must keep its original (default) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackingNumbers = @{
    2  = "320018612013"
    3  = "320018612024"
    4  = "320018612057"
    5  = "320018612079"
    6  = "320018612116"
    7  = "320018612138"
    8  = "320018612160"
    9  = "320018612182"
    10 = "320018612219"
    11 = "320018612230"
    12 = "320018612274"
    13 = "320018612296"
    14 = "320018612322"
    15 = "320018612344"
    16 = "320018612377"
    17 = "320018612399"
    18 = "320018612436"
    19 = "320018612458"
    20 = "320018612480"
    21 = "320018612506"
    22 = "320018612539"
    23 = "320018612540"
    24 = "320018612550"
    25 = "320018612561"
    26 = "320018612572"
}

foreach ($row in $trackingNumbers.Keys) {
    $cell = $ws.Range("P$row")
    # Force text storage so the numeric-looking tracking number isn't
    # coerced into a Number cell (matches the original t="s" cell type).
    $cell.NumberFormat = "@"
    $cell.Value = $trackingNumbers[$row]
    # Re-apply the sheet's default (un-styled) format so the temporary
    # text-number-format doesn't leave a stray style index behind -
    # Q-column on the same row carries the same default/unstyled format.
    $cell.Style = $ws.Range("Q$row").Style
}
